$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.885.11"
$ws.Range("E2").Value = "  -1.36%  "

# Row 3
$ws.Range("D3").Value = "2.454.94"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "579.22"
$ws.Range("E5").Value = "  -2.29%  "

# Row 6
$ws.Range("E6").Value = "  -4.37%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").Value = "0.512"
$ws.Range("E8").Value = "  -2.75%  "

# Row 9
$ws.Range("D9").Value = "2.456.16"
$ws.Range("E9").Value = "  -1.62%  "

# Row 10
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  -4.25%  "

# Row 11
$ws.Range("E11").Value = "  -1.18%  "

# Row 12
$ws.Range("E12").Value = "  -3.87%  "

# Row 13
$ws.Range("E13").Value = "  -3.26%  "

# Row 14
$ws.Range("D14").Value = "2.903.04"
$ws.Range("E14").Value = "  -1.65%  "

# Row 15
$ws.Range("D15").Value = "25.26"
$ws.Range("E15").Value = "  -3.85%  "

# Row 16
$ws.Range("D16").Value = "66.768.95"
$ws.Range("E16").Value = "  -1.27%  "

# Row 17
$ws.Range("E17").Value = "  -5.08%  "

# Row 18
$ws.Range("D18").Value = "2.455.23"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").Value = "11.31"
$ws.Range("E19").Value = "  -4.26%  "

# Row 20
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  -3.96%  "

# Row 21
$ws.Range("D21").Value = "352.92"
$ws.Range("E21").Value = "  -2.92%  "

# Row 22
$ws.Range("E22").Value = "  -2.40%  "

# Row 23
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").Value = "69.14"
$ws.Range("E24").Value = "  -2.85%  "

# Row 25
$ws.Range("D25").Value = "4.19"
$ws.Range("E25").Value = "  -7.98%  "

# Row 26
$ws.Range("E26").Value = "  -8.20%  "

# Row 27
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  -10.23%  "

# Row 28
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.38%  "

# Row 29
$ws.Range("D29").Value = "2.579.76"
$ws.Range("E29").Value = "  -1.84%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0891"
$ws.Range("E30").Value = "  -7.88%  "

# Row 31
$ws.Range("D31").Value = "504.64"
$ws.Range("E31").Value = "  -5.12%  "

# Row 32
$ws.Range("D32").Value = "7.74"
$ws.Range("E32").Value = "  -6.03%  "

# Row 33
$ws.Range("E33").Value = "  -5.73%  "

# Row 34
$ws.Range("E34").Value = "  -7.30%  "

# Row 35
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").Value = "158.03"
$ws.Range("E36").Value = "  -0.33%  "

# Row 37
$ws.Range("E37").Value = "  -9.29%  "

# Row 38
$ws.Range("D38").Value = "18.55"

# Row 39
$ws.Range("D39").Value = "18.39"
$ws.Range("E39").Value = "  -1.14%  "

# Row 40
$ws.Range("D40").Value = "1.34"
$ws.Range("E40").Value = "  -6.51%  "

# Row 41
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.325"
$ws.Range("E42").Value = "  -6.50%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.66"
$ws.Range("E43").Value = "  -6.99%  "

# Row 44
$ws.Range("E44").Value = "  -7.11%  "

# Row 45
$ws.Range("D45").Value = "38.76"
$ws.Range("E45").Value = "  -2.70%  "

# Row 46
$ws.Range("E46").Value = "  -7.97%  "

# Row 47
$ws.Range("D47").Value = "140.83"
$ws.Range("E47").Value = "  -3.32%  "

# Row 48
$ws.Range("D48").Value = "3.46"
$ws.Range("E48").Value = "  -6.19%  "

# Row 49
$ws.Range("D49").Value = "0.512"
$ws.Range("E49").Value = "  -6.60%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0251"
$ws.Range("E50").Value = "  -7.88%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0731"
$ws.Range("E51").Value = "  -2.51%  "
